# Update cryptos list values to match latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.426.04'
Set-TextValue $ws.Range('E2') '  +1.59%  '
Set-TextValue $ws.Range('D3') '1.863.35'
Set-TextValue $ws.Range('E3') '  +0.80%  '
Set-TextValue $ws.Range('E4') '  -0.11%  '
Set-TextValue $ws.Range('D5') '310.94'
Set-TextValue $ws.Range('E5') '  +0.35%  '
Set-TextValue $ws.Range('E6') '  -0.10%  '
Set-TextValue $ws.Range('D7') '0.4776'
Set-TextValue $ws.Range('E7') '  -0.15%  '
Set-TextValue $ws.Range('D8') '0.3760'
Set-TextValue $ws.Range('E8') '  +2.09%  '
Set-TextValue $ws.Range('D9') '0.07321'
Set-TextValue $ws.Range('E9') '  +1.19%  '
Set-TextValue $ws.Range('D10') '0.9336'
Set-TextValue $ws.Range('E10') '  +0.20%  '
Set-TextValue $ws.Range('D11') '20.63'
Set-TextValue $ws.Range('E11') '  +4.46%  '
Set-TextValue $ws.Range('D12') '0.07808'
Set-TextValue $ws.Range('E12') '  +1.18%  '
Set-TextValue $ws.Range('D13') '1.900.32'
Set-TextValue $ws.Range('E13') '  -0.51%  '
Set-TextValue $ws.Range('D14') '5.427'
Set-TextValue $ws.Range('E14') '  +1.59%  '
Set-TextValue $ws.Range('E15') '  +1.72%  '
Set-TextValue $ws.Range('D16') '90.38'
Set-TextValue $ws.Range('E16') '  +1.71%  '
Set-TextValue $ws.Range('D17') '1.012'
Set-TextValue $ws.Range('E17') '  -0.21%  '
Set-TextValue $ws.Range('D18') '0.000008871'
Set-TextValue $ws.Range('E18') '  +2.65%  '
Set-TextValue $ws.Range('E19') '  -0.15%  '
Set-TextValue $ws.Range('D20') '27.469.36'
Set-TextValue $ws.Range('E20') '  +1.63%  '
Set-TextValue $ws.Range('D21') '14.71'
Set-TextValue $ws.Range('E21') '  +1.42%  '
Set-TextValue $ws.Range('D22') '5.109'
Set-TextValue $ws.Range('E22') '  +1.01%  '
Set-TextValue $ws.Range('D23') '10.69'
Set-TextValue $ws.Range('E23') '  +0.34%  '
Set-TextValue $ws.Range('D24') '1.940'
Set-TextValue $ws.Range('E24') '  +0.54%  '
Set-TextValue $ws.Range('D25') '155.50'
Set-TextValue $ws.Range('E25') '  +1.75%  '
Set-TextValue $ws.Range('E26') '  +1.30%  '
Set-TextValue $ws.Range('D27') '2.017'
Set-TextValue $ws.Range('E27') '  +0.65%  '
Set-TextValue $ws.Range('D28') '115.34'
Set-TextValue $ws.Range('E28') '  +0.82%  '
Set-TextValue $ws.Range('D29') '4.940'
Set-TextValue $ws.Range('E29') '  -1.14%  '
Set-TextValue $ws.Range('D30') '0.08889'
Set-TextValue $ws.Range('E30') '  -0.18%  '
Set-TextValue $ws.Range('D31') '3.320'
Set-TextValue $ws.Range('E31') '  +0.85%  '
Set-TextValue $ws.Range('D32') '1.213'
Set-TextValue $ws.Range('E32') '  +3.10%  '
Set-TextValue $ws.Range('B33') 'ImmutableX'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D33') '0.7530'
Set-TextValue $ws.Range('E33') '  +1.00%  '
Set-TextValue $ws.Range('B34') 'Filecoin'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D34') '4.595'
Set-TextValue $ws.Range('E34') '  +1.90%  '
Set-TextValue $ws.Range('D35') '2.733'
Set-TextValue $ws.Range('E35') '  -0.35%  '
Set-TextValue $ws.Range('D36') '0.02048'
Set-TextValue $ws.Range('E36') '  +4.56%  '
Set-TextValue $ws.Range('D37') '1.116'
Set-TextValue $ws.Range('E37') '  +0.03%  '
Set-TextValue $ws.Range('D38') '0.05263'
Set-TextValue $ws.Range('E38') '  -0.06%  '
Set-TextValue $ws.Range('D39') '2.984'
Set-TextValue $ws.Range('E39') '  +0.17%  '
Set-TextValue $ws.Range('D40') '0.5323'
Set-TextValue $ws.Range('D41') '7.055'
Set-TextValue $ws.Range('E41') '  +0.74%  '
Set-TextValue $ws.Range('D42') '8.697'
Set-TextValue $ws.Range('E42') '  +5.86%  '
Set-TextValue $ws.Range('D43') '0.1525'
Set-TextValue $ws.Range('E43') '  +0.57%  '
Set-TextValue $ws.Range('D44') '10.60'
Set-TextValue $ws.Range('E44') '  +0.29%  '
Set-TextValue $ws.Range('D45') '0.4806'
Set-TextValue $ws.Range('E45') '  +1.10%  '
Set-TextValue $ws.Range('D46') '1.011'
Set-TextValue $ws.Range('E46') '  -0.14%  '
Set-TextValue $ws.Range('D47') '1.656'
Set-TextValue $ws.Range('E47') '  +2.61%  '
Set-TextValue $ws.Range('D48') '102.93'
Set-TextValue $ws.Range('E48') '  +0.78%  '
Set-TextValue $ws.Range('D49') '67.25'
Set-TextValue $ws.Range('E49') '  +2.53%  '
Set-TextValue $ws.Range('D50') '0.06076'
Set-TextValue $ws.Range('E50') '  +0.28%  '
Set-TextValue $ws.Range('D51') '0.9186'
Set-TextValue $ws.Range('E51') '  +3.39%  '
